$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 15 new rows at the top of the data (rows 1-15), shifting everything else down
$ws.Range("A1:A15").EntireRow.Insert()

# Populate the 15 new transaction rows
$ws.Range("A1").Value = 41823
$ws.Range("B1").Value = '  TRANSFERENCIA INTERNET'
$ws.Range("C1").Value = 'C'
$ws.Range("D1").Value = '0005026708'
$ws.Range("E1").Value = 'AG. NORTE'
$ws.Range("F1").Value = '475.18  '
$ws.Range("G1").Value = '3228.65'
$ws.Range("H1").Formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A1,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B1,"'', ''mo_tipo'' => ''",C1,"'', ''mo_documento'' => ''",D1,"'', ''mo_oficina'' => ''",E1,"'', ''mo_monto'' => ",F1,", ''mo_saldo'' => ",G1,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd H:m:s"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL, ''mo_fecha_borrado'' => NULL, ''mo_quien_borra'' => NULL, ''mo_borrado_logico'' => false),")'

$ws.Range("A2").Value = 41822
$ws.Range("B2").Value = 'INTERES A SU FAVOR'
$ws.Range("C2").Value = 'C'
$ws.Range("D2").Value = '0000825323'
$ws.Range("E2").Value = 'AGENCIA PARA PROCESOS BATCH'
$ws.Range("F2").Value = '0.17  '
$ws.Range("G2").Value = '2753.47'
$ws.Range("H2").Formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A2,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B2,"'', ''mo_tipo'' => ''",C2,"'', ''mo_documento'' => ''",D2,"'', ''mo_oficina'' => ''",E2,"'', ''mo_monto'' => ",F2,", ''mo_saldo'' => ",G2,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd H:m:s"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL, ''mo_fecha_borrado'' => NULL, ''mo_quien_borra'' => NULL, ''mo_borrado_logico'' => false),")'

$ws.Range("A3").Value = 41821
$ws.Range("B3").Value = 'INTERES A SU FAVOR'
$ws.Range("C3").Value = 'C'
$ws.Range("D3").Value = '0000825325'
$ws.Range("E3").Value = 'AGENCIA PARA PROCESOS BATCH'
$ws.Range("F3").Value = '0.17  '
$ws.Range("G3").Value = '2753.30'
$ws.Range("H3").Formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A3,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B3,"'', ''mo_tipo'' => ''",C3,"'', ''mo_documento'' => ''",D3,"'', ''mo_oficina'' => ''",E3,"'', ''mo_monto'' => ",F3,", ''mo_saldo'' => ",G3,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd H:m:s"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL, ''mo_fecha_borrado'' => NULL, ''mo_quien_borra'' => NULL, ''mo_borrado_logico'' => false),")'

$ws.Range("A4").Value = 41820
$ws.Range("B4").Value = 'INTERES A SU FAVOR'
$ws.Range("C4").Value = 'C'
$ws.Range("D4").Value = '0000825325'
$ws.Range("E4").Value = 'AGENCIA PARA PROCESOS BATCH'
$ws.Range("F4").Value = '0.17  '
$ws.Range("G4").Value = '2753.13'
$ws.Range("H4").Formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A4,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B4,"'', ''mo_tipo'' => ''",C4,"'', ''mo_documento'' => ''",D4,"'', ''mo_oficina'' => ''",E4,"'', ''mo_monto'' => ",F4,", ''mo_saldo'' => ",G4,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd H:m:s"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL, ''mo_fecha_borrado'' => NULL, ''mo_quien_borra'' => NULL, ''mo_borrado_logico'' => false),")'

$ws.Range("A5").Value = 41817
$ws.Range("B5").Value = 'INTERES A SU FAVOR'
$ws.Range("C5").Value = 'C'
$ws.Range("D5").Value = '0000825329'
$ws.Range("E5").Value = 'AGENCIA PARA PROCESOS BATCH'
$ws.Range("F5").Value = '0.52  '
$ws.Range("G5").Value = '2752.96'
$ws.Range("H5").Formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A5,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B5,"'', ''mo_tipo'' => ''",C5,"'', ''mo_documento'' => ''",D5,"'', ''mo_oficina'' => ''",E5,"'', ''mo_monto'' => ",F5,", ''mo_saldo'' => ",G5,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd H:m:s"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL, ''mo_fecha_borrado'' => NULL, ''mo_quien_borra'' => NULL, ''mo_borrado_logico'' => false),")'

$ws.Range("A6").Value = 41816
$ws.Range("B6").Value = 'INTERES A SU FAVOR'
$ws.Range("C6").Value = 'C'
$ws.Range("D6").Value = '0000825333'
$ws.Range("E6").Value = 'AGENCIA PARA PROCESOS BATCH'
$ws.Range("F6").Value = '0.17  '
$ws.Range("G6").Value = '2752.44'
$ws.Range("H6").Formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A6,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B6,"'', ''mo_tipo'' => ''",C6,"'', ''mo_documento'' => ''",D6,"'', ''mo_oficina'' => ''",E6,"'', ''mo_monto'' => ",F6,", ''mo_saldo'' => ",G6,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd H:m:s"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL, ''mo_fecha_borrado'' => NULL, ''mo_quien_borra'' => NULL, ''mo_borrado_logico'' => false),")'

$ws.Range("A7").Value = 41816
$ws.Range("B7").Value = '  TRANSFERENCIA INTERNET'
$ws.Range("C7").Value = 'C'
$ws.Range("D7").Value = '0003722472'
$ws.Range("E7").Value = 'AG. NORTE'
$ws.Range("F7").Value = '729.66  '
$ws.Range("G7").Value = '2752.27'
$ws.Range("H7").Formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A7,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B7,"'', ''mo_tipo'' => ''",C7,"'', ''mo_documento'' => ''",D7,"'', ''mo_oficina'' => ''",E7,"'', ''mo_monto'' => ",F7,", ''mo_saldo'' => ",G7,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd H:m:s"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL, ''mo_fecha_borrado'' => NULL, ''mo_quien_borra'' => NULL, ''mo_borrado_logico'' => false),")'

$ws.Range("A8").Value = 41815
$ws.Range("B8").Value = 'INTERES A SU FAVOR'
$ws.Range("C8").Value = 'C'
$ws.Range("D8").Value = '0000825339'
$ws.Range("E8").Value = 'AGENCIA PARA PROCESOS BATCH'
$ws.Range("F8").Value = '0.13  '
$ws.Range("G8").Value = '2022.61'
$ws.Range("H8").Formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A8,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B8,"'', ''mo_tipo'' => ''",C8,"'', ''mo_documento'' => ''",D8,"'', ''mo_oficina'' => ''",E8,"'', ''mo_monto'' => ",F8,", ''mo_saldo'' => ",G8,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd H:m:s"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL, ''mo_fecha_borrado'' => NULL, ''mo_quien_borra'' => NULL, ''mo_borrado_logico'' => false),")'

$ws.Range("A9").Value = 41814
$ws.Range("B9").Value = 'INTERES A SU FAVOR'
$ws.Range("C9").Value = 'C'
$ws.Range("D9").Value = '0000825339'
$ws.Range("E9").Value = 'AGENCIA PARA PROCESOS BATCH'
$ws.Range("F9").Value = '0.13  '
$ws.Range("G9").Value = '2022.48'
$ws.Range("H9").Formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A9,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B9,"'', ''mo_tipo'' => ''",C9,"'', ''mo_documento'' => ''",D9,"'', ''mo_oficina'' => ''",E9,"'', ''mo_monto'' => ",F9,", ''mo_saldo'' => ",G9,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd H:m:s"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL, ''mo_fecha_borrado'' => NULL, ''mo_quien_borra'' => NULL, ''mo_borrado_logico'' => false),")'

$ws.Range("A10").Value = 41813
$ws.Range("B10").Value = 'INTERES A SU FAVOR'
$ws.Range("C10").Value = 'C'
$ws.Range("D10").Value = '0000825340'
$ws.Range("E10").Value = 'AGENCIA PARA PROCESOS BATCH'
$ws.Range("F10").Value = '0.13  '
$ws.Range("G10").Value = '2022.35'
$ws.Range("H10").Formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A10,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B10,"'', ''mo_tipo'' => ''",C10,"'', ''mo_documento'' => ''",D10,"'', ''mo_oficina'' => ''",E10,"'', ''mo_monto'' => ",F10,", ''mo_saldo'' => ",G10,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd H:m:s"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL, ''mo_fecha_borrado'' => NULL, ''mo_quien_borra'' => NULL, ''mo_borrado_logico'' => false),")'

$ws.Range("A11").Value = 41813
$ws.Range("B11").Value = 'PAGO/RETIRO'
$ws.Range("C11").Value = 'D'
$ws.Range("D11").Value = '0010083738'
$ws.Range("E11").Value = 'EL GIRON'
$ws.Range("F11").Value = '124.38  '
$ws.Range("G11").Value = '2022.22'
$ws.Range("H11").Formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A11,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B11,"'', ''mo_tipo'' => ''",C11,"'', ''mo_documento'' => ''",D11,"'', ''mo_oficina'' => ''",E11,"'', ''mo_monto'' => ",F11,", ''mo_saldo'' => ",G11,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd H:m:s"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL, ''mo_fecha_borrado'' => NULL, ''mo_quien_borra'' => NULL, ''mo_borrado_logico'' => false),")'

$ws.Range("A12").Value = 41810
$ws.Range("B12").Value = 'INTERES A SU FAVOR'
$ws.Range("C12").Value = 'C'
$ws.Range("D12").Value = '0000825384'
$ws.Range("E12").Value = 'AGENCIA PARA PROCESOS BATCH'
$ws.Range("F12").Value = '0.40  '
$ws.Range("G12").Value = '2146.60'
$ws.Range("H12").Formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A12,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B12,"'', ''mo_tipo'' => ''",C12,"'', ''mo_documento'' => ''",D12,"'', ''mo_oficina'' => ''",E12,"'', ''mo_monto'' => ",F12,", ''mo_saldo'' => ",G12,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd H:m:s"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL, ''mo_fecha_borrado'' => NULL, ''mo_quien_borra'' => NULL, ''mo_borrado_logico'' => false),")'

$ws.Range("A13").Value = 41809
$ws.Range("B13").Value = 'INTERES A SU FAVOR'
$ws.Range("C13").Value = 'C'
$ws.Range("D13").Value = '0000825389'
$ws.Range("E13").Value = 'AGENCIA PARA PROCESOS BATCH'
$ws.Range("F13").Value = '0.13  '
$ws.Range("G13").Value = '2146.20'
$ws.Range("H13").Formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A13,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B13,"'', ''mo_tipo'' => ''",C13,"'', ''mo_documento'' => ''",D13,"'', ''mo_oficina'' => ''",E13,"'', ''mo_monto'' => ",F13,", ''mo_saldo'' => ",G13,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd H:m:s"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL, ''mo_fecha_borrado'' => NULL, ''mo_quien_borra'' => NULL, ''mo_borrado_logico'' => false),")'

$ws.Range("A14").Value = 41808
$ws.Range("B14").Value = 'INTERES A SU FAVOR'
$ws.Range("C14").Value = 'C'
$ws.Range("D14").Value = '0000825394'
$ws.Range("E14").Value = 'AGENCIA PARA PROCESOS BATCH'
$ws.Range("F14").Value = '0.13  '
$ws.Range("G14").Value = '2146.07'
$ws.Range("H14").Formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A14,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B14,"'', ''mo_tipo'' => ''",C14,"'', ''mo_documento'' => ''",D14,"'', ''mo_oficina'' => ''",E14,"'', ''mo_monto'' => ",F14,", ''mo_saldo'' => ",G14,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd H:m:s"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL, ''mo_fecha_borrado'' => NULL, ''mo_quien_borra'' => NULL, ''mo_borrado_logico'' => false),")'

$ws.Range("A15").Value = 41807
$ws.Range("B15").Value = 'INTERES A SU FAVOR'
$ws.Range("C15").Value = 'C'
$ws.Range("D15").Value = '0000825395'
$ws.Range("E15").Value = 'AGENCIA PARA PROCESOS BATCH'
$ws.Range("F15").Value = '0.13  '
$ws.Range("G15").Value = '2145.94'
$ws.Range("H15").Formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A15,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B15,"'', ''mo_tipo'' => ''",C15,"'', ''mo_documento'' => ''",D15,"'', ''mo_oficina'' => ''",E15,"'', ''mo_monto'' => ",F15,", ''mo_saldo'' => ",G15,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd H:m:s"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL, ''mo_fecha_borrado'' => NULL, ''mo_quien_borra'' => NULL, ''mo_borrado_logico'' => false),")'

# Update the view selection to reflect the new shared-formula range
$ws.Range("H1:H15").Select()

